$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new data rows (13 and 14) to the "Artfynd" sheet, matching the two
# rows appended in the source diff. Columns are addressed by letter so the
# mapping to the diff is easy to follow.
# ---------------------------------------------------------------------------

function Set-TextCell($sheet, $row, $col, $text) {
    # Force the value to be stored as literal text even when it looks like a
    # number or a date (Excel would otherwise silently convert "187" to the
    # number 187, or "2023-08-11" to a date serial). A leading apostrophe is
    # the normal Excel mechanism for "store this as text".
    $sheet.Cells.Item($row, $col).Value = "'" + $text
}

$rows = @(
    @{
        Row = 13
        A = 111845306; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        I = "187"; J = "plantor/tuvor"; K = "blomning"
        P = "Brinken, SV (mellrsta) (knärot), Vstm"
        Q = 564522; R = 6615783; S = 75
        T = "Västmanland"; U = "Surahammar"; V = "Västmanland"; W = "Sura"
        X = "U-Sur-0535"
        Y = "2023-08-11"; AA = "2023-08-11"
        AC = "X: (7) 30 pl, A: 186/320 (5), 7 pl, 1 bl, B. 183/292 (4), 100 pl, 4 bl, C:205/297 (0)50 pl (5 m från tänkt basväg)"
        AD = $false; AE = $false; AG = $false
        AW = "Bo Eriksson"; AX = "Michael Lander"; AY = "Floraväkteri Sverige"
    },
    @{
        Row = 14
        A = 111845440; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        I = "230"; J = "plantor/tuvor"; K = "blomning"
        P = "Brinken, SV, (mellersta, syd) (knärot), Vstm"
        Q = 564551; R = 6615754; S = 25
        T = "Västmanland"; U = "Surahammar"; V = "Västmanland"; W = "Sura"
        X = "U-Sur-0536"
        Y = "2023-08-11"; AA = "2023-08-11"
        AC = "X. (7), 150 pl, A: 192/262 (4), 20 pl, 3 bl, (5 m från tänkt basväg B: 2307243 (4), 60 pl, 10 bl."
        AD = $false; AE = $false; AG = $false
        AW = "Bo Eriksson"; AX = "Michael Lander"; AY = "Floraväkteri Sverige"
    }
)

foreach ($r in $rows) {
    $row = $r.Row

    # --- plain numeric columns -------------------------------------------------
    $ws.Cells.Item($row, 1).Value  = $r.A    # A  Id
    $ws.Cells.Item($row, 2).Value  = $r.B    # B  Taxonsorteringsordning
    $ws.Cells.Item($row, 5).Value  = $r.E    # E  TaxonId
    $ws.Cells.Item($row, 17).Value = $r.Q    # Q  Ost
    $ws.Cells.Item($row, 18).Value = $r.R    # R  Nord
    $ws.Cells.Item($row, 19).Value = $r.S    # S  Noggrannhet

    # --- plain text columns -----------------------------------------------------
    $ws.Cells.Item($row, 3).Value  = $r.C    # C  Valideringsstatus
    $ws.Cells.Item($row, 4).Value  = $r.D    # D  Rödlistade
    $ws.Cells.Item($row, 6).Value  = $r.F    # F  Artnamn
    $ws.Cells.Item($row, 7).Value  = $r.G    # G  Vetenskapligt namn
    $ws.Cells.Item($row, 8).Value  = $r.H    # H  Auktor
    $ws.Cells.Item($row, 10).Value = $r.J    # J  Enhet
    $ws.Cells.Item($row, 11).Value = $r.K    # K  Ålder-Stadium
    $ws.Cells.Item($row, 16).Value = $r.P    # P  Lokalnamn
    $ws.Cells.Item($row, 20).Value = $r.T    # T  Län
    $ws.Cells.Item($row, 21).Value = $r.U    # U  Kommun
    $ws.Cells.Item($row, 22).Value = $r.V    # V  Provins
    $ws.Cells.Item($row, 23).Value = $r.W    # W  Församling
    $ws.Cells.Item($row, 24).Value = $r.X    # X  Externid
    $ws.Cells.Item($row, 29).Value = $r.AC   # AC Publik kommentar
    $ws.Cells.Item($row, 49).Value = $r.AW   # AW Rapportör
    $ws.Cells.Item($row, 50).Value = $r.AX   # AX Observatörer
    $ws.Cells.Item($row, 51).Value = $r.AY   # AY Projektnamn

    # --- text columns that must not be auto-converted to number/date ------------
    Set-TextCell $ws $row 9  $r.I     # I  Antal (stored as text, e.g. "187")
    Set-TextCell $ws $row 25 $r.Y     # Y  Startdatum
    Set-TextCell $ws $row 27 $r.AA    # AA Slutdatum

    # --- boolean columns ----------------------------------------------------
    $ws.Cells.Item($row, 30).Value = $r.AD   # AD Ej återfunnen
    $ws.Cells.Item($row, 31).Value = $r.AE   # AE Osäker artbestämning
    $ws.Cells.Item($row, 33).Value = $r.AG   # AG Ospontan
}
